# ecYaliGEM curated for more accurate acetaldehyde metabolism
# Adds 5 new curated kcat entries (rows 38-42) to the "customKcats" sheet.
#
# Cell writes are ordered to reproduce the author's original shared-string
# insertion order (reactions/genes/proteins entered first, notes last).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customKcats")

$note = 'Putative enzyme. Using the lowest kcat instead, which is the one predicted by DLKcat in this case.'

# Row 38 : y000173 / YALI0C03025g / Q6CD79
$ws.Range("E38").Value = 'y000173'
$ws.Range("B38").Value = 'YALI0C03025g'
$ws.Range("C38").Value = 'YALI0C03025g'
$ws.Range("A38").Value = 'Q6CD79'
$ws.Range("F38").Value = $note

# Row 39 : y000185 / YALI0D07942g / Q6C9V7
$ws.Range("B39").Value = 'YALI0D07942g'
$ws.Range("C39").Value = 'YALI0D07942g'
$ws.Range("A39").Value = 'Q6C9V7'
$ws.Range("E39").Value = 'y000185'

# Row 40 : y000185 / YALI0F04444g / Q6C2W9
$ws.Range("B40").Value = 'YALI0F04444g'
$ws.Range("C40").Value = 'YALI0F04444g'
$ws.Range("A40").Value = 'Q6C2W9'

# Row 41 : y002116 / YALI0D07942g / Q6C9V7
$ws.Range("E41").Value = 'y002116'
$ws.Range("F39").Value = $note
$ws.Range("E40").Value = 'y000185'
$ws.Range("F40").Value = $note
$ws.Range("A41").Value = 'Q6C9V7'
$ws.Range("B41").Value = 'YALI0D07942g'
$ws.Range("C41").Value = 'YALI0D07942g'
$ws.Range("F41").Value = $note

# Row 42 : y002116 / YALI0F04444g / Q6C2W9
$ws.Range("A42").Value = 'Q6C2W9'
$ws.Range("B42").Value = 'YALI0F04444g'
$ws.Range("C42").Value = 'YALI0F04444g'
$ws.Range("E42").Value = 'y002116'
$ws.Range("F42").Value = $note

# kcat values (column D) and stoichiometry (column G)
$ws.Range("D38").Value = 12.3903
$ws.Range("G38").Value = 1
$ws.Range("D39").Value = 26.5155
$ws.Range("G39").Value = 1
$ws.Range("D40").Value = 32.9861
$ws.Range("G40").Value = 1
$ws.Range("D41").Value = 6.7826
$ws.Range("G41").Value = 1
$ws.Range("D42").Value = 9.1103
$ws.Range("G42").Value = 1

# Final cursor / selection position left by the editing session.
$ws.Range("F44").Select() | Out-Null
